$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between the paired rows below (match
# scheduling details were swapped between two fixtures that share the same
# match date). Columns A,B,C,D,E,J keep their original row position.
$swapCols = @("F","G","H","I","K","L","M","N","O","P","Q","R","S","T","U","V")

function Swap-RowData($row1, $row2) {
    foreach ($col in $swapCols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Rows 7 & 8 swap (TS Galaxy vs Cape Town Spurs <-> Kaizer Chiefs vs Chippa Utd.)
Swap-RowData 7 8

# Rows 86 & 87 swap (Golden Arrows vs AmaZulu <-> Richards Bay vs Cape Town Spurs)
Swap-RowData 86 87

# Rows 92 & 93 swap (TS Galaxy vs Polokwane <-> Swallows vs Kaizer Chiefs)
Swap-RowData 92 93

# Append new row 96 with the latest fixture result.
$ws.Range("A95:V95").Copy()
$ws.Range("A96:V96").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A96").Value = 95
$ws.Range("B96").Value = "south-africa"
$ws.Range("C96").Value = "premier-league"
$ws.Range("D96").Value = "2023-2024"
$ws.Range("E96").Value = 45259.77083333334
$ws.Range("F96").Value = "Supersport Utd"
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = "Mamelodi Sundowns"
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = 3.72
$ws.Range("K96").Value = "22/11/2023 18:43"
$ws.Range("L96").Value = 4.9
$ws.Range("M96").Value = "29/11/2023 18:26"
$ws.Range("N96").Value = 3.04
$ws.Range("O96").Value = "22/11/2023 18:43"
$ws.Range("P96").Value = 3.24
$ws.Range("Q96").Value = "29/11/2023 18:21"
$ws.Range("R96").Value = 2.19
$ws.Range("S96").Value = "22/11/2023 18:43"
$ws.Range("T96").Value = 1.87
$ws.Range("U96").Value = "29/11/2023 18:26"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/south-africa/premier-league/supersport-utd-mamelodi-sundowns/EazKsdKf/"
